$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header cell A1 changes from numeric 0 to text "Description"
$ws.Range("A1").Value = "Description"

# Row 2
$ws.Range("A2").Value = "Move Robot21 to location (6, 12) and remove the toolkit."
$ws.Range("B2").Value = $false
$ws.Range("F2").Value = $false

# Row 3
$ws.Range("A3").Value = "Move Robot41 to location (1, 11) and remove the liquid spill."

# Row 4
$ws.Range("A4").Value = "Move Robot9 to location (1, 5) and remove the large debris."

# Row 5
$ws.Range("A5").Value = "Move Robot42 to location (1, 11) and remove the dust."
$ws.Range("B5").Value = $true

# Row 6
$ws.Range("A6").Value = "Move Robot32 to location (3, 9) and remove the grass."
$ws.Range("B6").Value = $false
$ws.Range("F6").Value = $false

# Row 7
$ws.Range("A7").Value = "Move Robot14 to location (11, 12) and remove the small debris."

# Row 8
$ws.Range("A8").Value = "Move Robot39 to location (6, 4) and remove the vehicle."
$ws.Range("B8").Value = $false
$ws.Range("F8").Value = $false

# Row 9
$ws.Range("A9").Value = "Move Robot15 to location (11, 2) and remove the construction materials."

# Row 10
$ws.Range("A10").Value = "Move Robot2 to location (1, 10) and remove the tree branches."

# Row 11
$ws.Range("A11").Value = "Move Robot26 to location (1, 3) and remove the screws."
$ws.Range("B11").Value = $false
$ws.Range("F11").Value = $false
